$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.727.62'
$ws.Range('E2').Value = '  -0.06%  '
$ws.Range('D3').Value = '2.076.40'
$ws.Range('E3').Value = '  -1.23%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = '233.50'
$ws.Range('E5').Value = '  -0.57%  '
$ws.Range('D6').Value = '0.623'
$ws.Range('E6').Value = '  -0.11%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').Value = '58.11'
$ws.Range('E8').Value = '  +0.27%  '
$ws.Range('D9').Value = '0.395'
$ws.Range('E9').Value = '  +0.95%  '
$ws.Range('D10').Value = '0.0783'
$ws.Range('E10').Value = '  +0.68%  '
$ws.Range('D11').Value = '0.106'
$ws.Range('E11').Value = '  +2.61%  '
$ws.Range('D12').Value = '2.380.43'
$ws.Range('E12').Value = '  -1.44%  '
$ws.Range('D13').Value = '14.75'
$ws.Range('E13').Value = '  +1.51%  '
$ws.Range('D14').Value = '20.85'
$ws.Range('E14').Value = '  -2.12%  '
$ws.Range('D15').Value = '0.774'
$ws.Range('E15').Value = '  -0.95%  '
$ws.Range('D16').Value = '5.32'
$ws.Range('E16').Value = '  +2.04%  '
$ws.Range('D17').Value = '2.155.21'
$ws.Range('E17').Value = '  +2.24%  '
$ws.Range('D18').Value = '37.652.81'
$ws.Range('E18').Value = '  -0.21%  '
$ws.Range('E19').Value = '  -1.05%  '
$ws.Range('E20').Value = '  +1.26%  '
$ws.Range('E21').Value = '  +1.24%  '
$ws.Range('D22').Value = '228.13'
$ws.Range('E23').Value = '  -0.04%  '
$ws.Range('E24').Value = '  -0.88%  '
$ws.Range('D25').Value = '2.39'
$ws.Range('E25').Value = '  -0.59%  '
$ws.Range('D26').Value = '170.66'
$ws.Range('E26').Value = '  +1.52%  '
$ws.Range('E27').Value = '  +5.38%  '
$ws.Range('E28').Value = '  +1.21%  '
$ws.Range('E29').Value = '  -0.10%  '
$ws.Range('E30').Value = '  -2.25%  '
$ws.Range('E31').Value = '  +2.36%  '
$ws.Range('D32').Value = '4.68'
$ws.Range('E32').Value = '  +1.21%  '
$ws.Range('E33').Value = '  +1.12%  '
$ws.Range('D34').Value = '4.64'
$ws.Range('E34').Value = '  +1.46%  '
$ws.Range('E35').Value = '  -4.05%  '
$ws.Range('E36').Value = '  +1.40%  '
$ws.Range('D37').Value = '3.39'
$ws.Range('E37').Value = '  -2.07%  '
$ws.Range('E38').Value = '  -0.07%  '
$ws.Range('E39').Value = '  -2.52%  '
$ws.Range('D40').Value = '0.0968'
$ws.Range('E40').Value = '  -3.05%  '
$ws.Range('D41').Value = '98.26'
$ws.Range('E41').Value = '  +1.88%  '
$ws.Range('E42').Value = '  -2.32%  '
$ws.Range('D43').Value = '0.0215'
$ws.Range('E43').Value = '  +0.73%  '
$ws.Range('D44').Value = '1.448.66'
$ws.Range('E44').Value = '  -1.77%  '
$ws.Range('B45').Value = 'InjectiveProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D45').Value = '16.37'
$ws.Range('E45').Value = '  +6.49%  '
$ws.Range('B46').Value = 'TrustWalletToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D46').Value = '1.15'
$ws.Range('E46').Value = '  -1.86%  '
$ws.Range('D47').Value = '4.18'
$ws.Range('E47').Value = '  +1.64%  '
$ws.Range('E48').Value = '  +0.86%  '
$ws.Range('D49').Value = '7.40'
$ws.Range('E49').Value = '  +1.16%  '
$ws.Range('D50').Value = '3.01'
$ws.Range('E50').Value = '  -0.55%  '
$ws.Range('D51').Value = '2.265.32'
$ws.Range('E51').Value = '  -1.51%  '
